$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "FAIL"
}

$ws.Range("B2:B11").Select()
$ws.Range("B3").Activate()
